$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the legend emoji / color-name labels used in the "statut" columns.
# Mapping (synthetic array refactor): black -> blue, keep others positionally aligned.
$ws.Cells.Replace("⬛", "📘")
$ws.Cells.Replace("🟥", "📕")
$ws.Cells.Replace("🟩", "📗")
$ws.Cells.Replace("🟧", "📙")
$ws.Cells.Replace("noir", "bleu")
